$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.097.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.232.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.228.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.761.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.227.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.064.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0741"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0397"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "410.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.864.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  +0.43%  "
